$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37; everything from row 37 downward shifts to row+1.
$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the new data point (template columns copied
# from the surrounding rows; only D/L/M/N/O/P/S vary row to row).
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 44518
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100107
$ws.Range("H37").Value = "Otros"
$ws.Range("I37").Value = 100107002
$ws.Range("J37").Value = "Chirimoya"
$ws.Range("K37").Value = "Cultivar IV Región"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 210
$ws.Range("N37").Value = 20000
$ws.Range("O37").Value = 20000
$ws.Range("P37").Value = 20000
$ws.Range("Q37").Value = "$/bandeja 10 kilos"
$ws.Range("R37").Value = "Provincia de Limarí"
$ws.Range("S37").Value = 2000
$ws.Range("T37").Value = 10

# Make sure the date cell keeps the same date number format as the rest of
# column D.
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
